$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its original location (end of the
#    "...en un formato intermedio." bullet). It is re-created further
#    down, inside the restructured Testabilidad bullet.
# ---------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Restructure the "Testabilidad: ..." bullet:
#      - the bullet itself keeps only the bold "Testabilidad: " label
#        (split into "Testabilidad" + ": " runs);
#      - the old body sentence becomes its own sub-bullet
#        ("Debe conseguirse un sistema sencillo...");
#      - a brand-new sub-bullet is appended ("Debe garantizarse que la
#        conversión al formato empleado...") with "al" underlined and
#        the relocated _GoBack bookmark right after "empleado".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Testabilidad: Debe conseguirse*correctas.", $false, $false, $true, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$titlePara = $rng.Paragraphs(1)
$paraEnd = $titlePara.Range.End
$titleLen = 14   # Len("Testabilidad: ")

# Text that currently follows the "Testabilidad: " label.
$bodyText = $d.Range($titlePara.Range.Start + $titleLen, $paraEnd).Text

# Two fresh paragraphs after the label bullet (they inherit its pPr/numPr).
$titlePara.Range.InsertParagraphAfter()
$titlePara.Range.InsertParagraphAfter()
$firstSubPara = $titlePara.Next()
$secondSubPara = $firstSubPara.Next()

# Move the old body sentence into the first new sub-bullet...
$firstSubPara.Range.Text = $bodyText
# ...and strip it back out of the label paragraph, leaving "Testabilidad: ".
$d.Range($titlePara.Range.Start + $titleLen, $paraEnd).Text = ""

# Demote both new bullets to the second list level (w:ilvl = 1).
$firstSubPara.Range.ListFormat.ListLevelNumber = 2
$secondSubPara.Range.ListFormat.ListLevelNumber = 2

# Fill in the brand-new sentence about format conversion.
$newSentence = "Debe garantizarse que la conversión al formato empleado por la aplicación " + `
                "es correcta y no se pierde la información en el proceso."
$secondSubPara.Range.Text = $newSentence

# Underline "al" (as two single-letter runs, matching how Word itself
# would split them while the user was editing).
$alStart = $secondSubPara.Range.Start + 36
$d.Range($alStart, $alStart + 1).Font.Underline = 1
$d.Range($alStart + 1, $alStart + 2).Font.Underline = 1

# Isolate "empleado" into its own run so the relocated bookmark sits
# cleanly between it and the following " por la aplicación..." text.
$empStart = $secondSubPara.Range.Start + 47
$empEnd = $secondSubPara.Range.Start + 55
$empRange = $d.Range($empStart, $empEnd)
$empRange.Bold = 1
$empRange.Bold = 0

# Re-create the _GoBack bookmark right after "empleado".
$d.Bookmarks.Add("_GoBack", $d.Range($empEnd, $empEnd))

# Split the bold label "Testabilidad: " into "Testabilidad" + ": " runs
# (mirrors the proofing-mark boundary Word itself introduces there).
$splitPoint = $titlePara.Range.Start + 12
$colonRange = $d.Range($splitPoint, $splitPoint + 2)
$colonRange.Bold = 0
$colonRange.Bold = 1

Write-Output "Testabilidad section restructured"
